$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the values previously in row 3 (Mesa = C2 stays 1)
$ws.Range("A2").Value = "Paolo"

# Force B2 to stay text (otherwise "+549" auto-converts to the number 549)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "+549"

# Remove row 3 entirely (used range shrinks to A1:R2)
$ws.Rows(3).Delete()
